$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Tabelle1" to "BFeld"
$ws.Name = "BFeld"

# Add new row of data (A15 = 112, B15 = 105)
$ws.Range("A15").Value = 112
$ws.Range("B15").Value = 105

# Update the selected cell to B15
$ws.Range("B15").Select()
